# [ADD] slides session 20
#
# The "Recommender Systems - Hands-On" slide instructs students to open a
# notebook file. Rename the referenced notebook from
#   recommender_system_exercises.ipynb
# to
#   part1_recommendation_system_exercises.ipynb
# wherever it is mentioned in the deck, preserving the surrounding text
# and run formatting.

$oldSnippet = ': "recommender_system_exercises.ipynb".'
$newSnippet = ': "part1_recommendation_system_exercises.ipynb".'

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $shp = $s.Shapes.Item($shi)
        if (-not $shp.HasTextFrame) { continue }

        $tr = $shp.TextFrame.TextRange
        $fullText = $tr.Text
        $startPos = $fullText.IndexOf($oldSnippet)
        while ($startPos -ge 0) {
            $sub = $tr.Characters($startPos + 1, $oldSnippet.Length)
            $sub.Text = $newSnippet

            $fullText = $tr.Text
            $startPos = $fullText.IndexOf($oldSnippet)
        }
    }
}
